# Update marksheet correct/total marks on the "quiz" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 "Marking" -> Right column value: 3 -> 5
$ws.Range("B11").Value = 5

# Row 12 "Total" -> Right column value: 81 -> 135
$ws.Range("B12").Value = 135

# Row 12 "Total" -> Max column (E12) text: 80/84 -> 135/140
$ws.Range("E12").Value = "135/140"
